$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("B2").Value = 7179
$ws.Range("C3").Value = 159935
$ws.Range("C4").Value = 150971
$ws.Range("C7").Value = 5.6
$ws.Range("C8").Value = 64.27
